# Generate Report for Handoff
# The "c9f9b0cf-cf70-4f22-89d1-741540a3278d.md" file has moved from
# "In Translation" to "Ready for handoff" with a fresh handoff timestamp
# and a priority change from "ht" to "mt". Reflect this on all three
# sheets (Overview summary + the per-language zh-cn / de-de detail
# sheets), then autofit the columns whose text just got longer.

$wb = $excel.ActiveWorkbook

# --- Overview sheet (per-file summary) ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-18 22:13:36"
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()

# --- zh-cn detail sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("H3").Value = "2016-08-18 22:13:31"
$ws.Columns.Item(3).AutoFit()

# --- de-de detail sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("H3").Value = "2016-08-18 22:13:36"
$ws.Columns.Item(3).AutoFit()
